$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 112. This pushes the
# existing rows 112-120 down to 114-122 (their contents are unchanged by
# the shift, matching the diff's row 114..122 == old row 112..120 data).
$ws.Rows.Item(112).Insert()
$ws.Rows.Item(112).Insert()

# New row 112 (weekly "Fruta" entry): Murcott / Primera, $/caja 10 kilos
$ws.Range("A112").Value = 4
$ws.Range("B112").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C112").Value = "Los Lagos"
$ws.Range("D112").Value = 44491
$ws.Range("E112").Value = 10
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100102
$ws.Range("H112").Value = "Cítricos"
$ws.Range("I112").Value = 100102004
$ws.Range("J112").Value = "Mandarina"
$ws.Range("K112").Value = "Murcott"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 700
$ws.Range("N112").Value = 5500
$ws.Range("O112").Value = 6000
$ws.Range("P112").Value = 5750
$ws.Range("Q112").Value = "$/caja 10 kilos"
$ws.Range("R112").Value = "Provincia de Limarí"
$ws.Range("S112").Value = 575
$ws.Range("T112").Value = 10

# New row 113 (weekly "Fruta" entry): Murcott / Segunda, $/caja 10 kilos
$ws.Range("A113").Value = 4
$ws.Range("B113").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C113").Value = "Los Lagos"
$ws.Range("D113").Value = 44491
$ws.Range("E113").Value = 10
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100102
$ws.Range("H113").Value = "Cítricos"
$ws.Range("I113").Value = 100102004
$ws.Range("J113").Value = "Mandarina"
$ws.Range("K113").Value = "Murcott"
$ws.Range("L113").Value = "Segunda"
$ws.Range("M113").Value = 300
$ws.Range("N113").Value = 4000
$ws.Range("O113").Value = 4000
$ws.Range("P113").Value = 4000
$ws.Range("Q113").Value = "$/caja 10 kilos"
$ws.Range("R113").Value = "Provincia de Limarí"
$ws.Range("S113").Value = 400
$ws.Range("T113").Value = 10
